# MAJ des tests faits
# Fill in the "Suivi des tests" table on the first sheet (Feuil1).
# The values are entered in the same order the original author typed
# them in Excel so that the shared-strings table is built up in the
# same sequence as in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: premier test ---
$ws.Range("A2").Value = "post avec un ticket test"
$ws.Range("B2").Value = "réussi"

# --- Row 3: résultat ---
$ws.Range("B3").Value = "echec"

# --- Row 4: résultat ---
$ws.Range("B4").Value = "réussite après modif: validation des donnes et types"

# --- Row 5: objet ---
$ws.Range("A5").Value = "sur le site (avant cypress)"

# --- Row 5: causes possibles ---
$ws.Range("C5").Value = "pas de methode fetch"

# --- Objets (col A) for rows 3 and 4 ---
$ws.Range("A3").Value = "Bdd sqlite 1"
$ws.Range("A4").Value = "Bdd sqlite 2"

# --- Header row ---
$ws.Range("A1").Value = "objet"
$ws.Range("B1").Value = "résultat"
$ws.Range("C1").Value = "causes possibles"

# --- Causes possibles (col C) for row 3 ---
$ws.Range("C3").Value = "erreur dans le code ?"

# --- Row 5 résultat (reuses "echec" already used in row 3) ---
$ws.Range("B5").Value = "echec"

# Auto-fit the two columns that contain the longest text so the
# columns are readable, like the original author did.
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Click away from the table, ending the selection on A16, as in the
# final saved state of the workbook.
$ws.Range("A16").Select()
